# Zeitplan.xlsx update: expand schedule with two shoot-day date headers and
# the associated new/shifted time-slot rows.
#
# -4122 == xlPasteFormats (PasteSpecial paste type). Using Copy + PasteSpecial
# lets Excel reuse/record the exact same cellXf (style) index it would when a
# user does this interactively, instead of us trying to hand-roll Font/Interior
# property writes that might resolve to a different (but visually equivalent)
# style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room -----------------------------------------------------
# A new "day" header row is inserted above the existing row 2, pushing the
# whole schedule (old rows 2-18) down by one (-> new rows 3-19).
$ws.Rows.Item(2).Insert()

# The gap between the two schedule blocks (old rows 19-23, never populated)
# shrinks from 5 rows down to a single new "day" header row. Remove the
# 4 surplus (already-empty) rows so the second block lands on new row 21
# with the header on new row 20.
$ws.Range("A20:A23").EntireRow.Delete()

# --- 2. Day headers -----------------------------------------------------
# New row 2: date header for the first shoot day (2025-01-14), bold font +
# date number format, same as the other header-ish cells (based on A1's
# bold style, numFmtId 14).
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = 45671

# New row 20: date header for the second shoot day (2025-01-17). Same style
# as the row-2 header above.
$ws.Range("A2").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 45674

# --- 3. Re-point the timestamps within the first day's block ------------
# The old "11:50" timestamp cell (previously on row 8 / now row 9) is
# cleared -- the timestamp now lives a few rows further down (new row 13),
# tagged with the "ACTOR#1 walks to the media-tech room" block instead.
$ws.Range("A9").ClearContents()

# New row 13 gets the 11:50 timestamp, with the highlighted (fillId 5)
# style used for the other "moving between rooms" rows in that block.
$ws.Range("A7").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").NumberFormat = "h:mm"
$ws.Range("A13").Value = 0.49305555555555558

# New row 14 (13:00) and new row 19 (13:40) are additional timestamps
# introduced partway through the first day's block, both using the style
# already used for the 11:50 slot before it moved (fillId 2).
$ws.Range("A9").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 0.54166666666666663

$ws.Range("A9").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = 0.56944444444444442

# --- 4. Timestamps for the second day's block ---------------------------
# New row 21 (13:00), matching style of the block's first row (fillId 3).
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = 0.54166666666666663

# New row 25 (13:40), matching the fillId-5 style used for the 11:50 slot.
$ws.Range("A13").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = 0.56944444444444442

# --- 5. Leave the selection where the author left off --------------------
[void]$ws.Range("B27").Select()
